$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 34 (shifts existing rows 34:39 down to 35:40,
# carrying the existing row-34 formatting onto the new blank row).
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the average air-pressure figure.
$ws.Range("A34").Value = "Ave pressure kpa (air)"
$ws.Range("B34").Value = 64.48

# The fill-only style that used to be applied to B3:B5 is no longer used;
# drop it so those cells fall back to the default "Normal" style.
$ws.Range("B3:B5").Style = "Normal"

# Reflect where the user had scrolled/selected when the workbook was saved.
[void]$excel.ActiveWindow.ScrollRow
$excel.ActiveWindow.ScrollRow = 6
[void]$ws.Range("F33").Select()
